$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = new price text (optional); E = new volume text }
$updates = @(
    @{ Row = 2; D = "90.660.73"; E = "  +0.16%  " }
    @{ Row = 3; D = "3.109.44"; E = "  +0.30%  " }
    @{ Row = 4; D = $null; E = "  +0.14%  " }
    @{ Row = 5; D = "242.74"; E = "  +2.91%  " }
    @{ Row = 6; D = "623.64"; E = "  -0.31%  " }
    @{ Row = 7; D = "1.18"; E = "  +14.84%  " }
    @{ Row = 8; D = "0.370"; E = "  +5.43%  " }
    @{ Row = 9; D = $null; E = "  -0.02%  " }
    @{ Row = 10; D = "3.106.53"; E = "  -8.14%  " }
    @{ Row = 11; D = "0.758"; E = "  +6.39%  " }
    @{ Row = 12; D = $null; E = "  +3.89%  " }
    @{ Row = 13; D = "0.0000251"; E = "  +3.90%  " }
    @{ Row = 14; D = "35.31"; E = "  -2.03%  " }
    @{ Row = 15; D = $null; E = "  -1.36%  " }
    @{ Row = 16; D = "90.509.84"; E = "  +0.23%  " }
    @{ Row = 17; D = "3.685.61"; E = "  +0.08%  " }
    @{ Row = 18; D = "3.111.67"; E = "  +1.51%  " }
    @{ Row = 19; D = "3.79"; E = "  +2.71%  " }
    @{ Row = 20; D = "14.35"; E = "  +0.74%  " }
    @{ Row = 21; D = $null; E = "  -0.17%  " }
    @{ Row = 22; D = "5.81"; E = "  +8.33%  " }
    @{ Row = 23; D = "448.63"; E = "  +0.28%  " }
    @{ Row = 24; D = "9.07"; E = "  +2.03%  " }
    @{ Row = 25; D = $null; E = "  -3.21%  " }
    @{ Row = 26; D = "93.64"; E = "  +4.29%  " }
    @{ Row = 27; D = "11.91"; E = "  -2.35%  " }
    @{ Row = 28; D = $null; E = "  +0.11%  " }
    @{ Row = 29; D = $null; E = "  +0.05%  " }
    @{ Row = 30; D = "0.176"; E = "  +11.35%  " }
    @{ Row = 31; D = "0.228"; E = "  +15.70%  " }
    @{ Row = 32; D = "9.11"; E = "  -1.60%  " }
    @{ Row = 33; D = $null; E = "  +35.58%  " }
    @{ Row = 34; D = "0.113"; E = "  +32.08%  " }
    @{ Row = 35; D = $null; E = "  +5.45%  " }
    @{ Row = 36; D = "26.63"; E = "  -2.43%  " }
    @{ Row = 37; D = "7.63"; E = "  +9.30%  " }
    @{ Row = 38; D = "4.18"; E = "  +29.36%  " }
    @{ Row = 39; D = $null; E = "  -0.02%  " }
    @{ Row = 40; D = "493.56"; E = "  -2.35%  " }
    @{ Row = 41; D = "3.60"; E = "  -4.25%  " }
    @{ Row = 42; D = $null; E = "  -1.76%  " }
    @{ Row = 43; D = $null; E = "  -0.81%  " }
    @{ Row = 44; D = $null; E = "  -0.30%  " }
    @{ Row = 46; D = "157.37"; E = "  +5.78%  " }
    @{ Row = 47; D = $null; E = "  -2.35%  " }
    @{ Row = 48; D = $null; E = "  -0.98%  " }
    @{ Row = 49; D = "4.58"; E = "  +1.61%  " }
    @{ Row = 50; D = "44.99"; E = "  +0.34%  " }
    @{ Row = 51; D = $null; E = "  +0.08%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force text storage so values like "0.370" / "3.60" keep their
        # exact printed form instead of being coerced to numbers.
        $ws.Range("D$($u.Row)").NumberFormat = "@"
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
